$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple per-cell updates (Price / Volume columns) for rows that did not reorder.
# ForceText=$true marks values that Excel would otherwise silently reinterpret as
# numbers (e.g. "152.40" -> 152.4), stripping the significant trailing zero.
$updates = @(
    @{Cell='D2'; Value='60.949.35'; ForceText=$false},
    @{Cell='E2'; Value='  -1.55%  '; ForceText=$false},
    @{Cell='D3'; Value='2.413.80'; ForceText=$false},
    @{Cell='E3'; Value='  -1.25%  '; ForceText=$false},
    @{Cell='E4'; Value='  -0.07%  '; ForceText=$false},
    @{Cell='D5'; Value='567.57'; ForceText=$false},
    @{Cell='E5'; Value='  -2.23%  '; ForceText=$false},
    @{Cell='D6'; Value='138.85'; ForceText=$false},
    @{Cell='E6'; Value='  -2.06%  '; ForceText=$false},
    @{Cell='E7'; Value='  +0.26%  '; ForceText=$false},
    @{Cell='D8'; Value='0.536'; ForceText=$false},
    @{Cell='E8'; Value='  +1.13%  '; ForceText=$false},
    @{Cell='D9'; Value='2.397.32'; ForceText=$false},
    @{Cell='E9'; Value='  -1.68%  '; ForceText=$false},
    @{Cell='E10'; Value='  -3.80%  '; ForceText=$false},
    @{Cell='E11'; Value='  -0.34%  '; ForceText=$false},
    @{Cell='D12'; Value='5.04'; ForceText=$false},
    @{Cell='E12'; Value='  -2.73%  '; ForceText=$false},
    @{Cell='D13'; Value='0.336'; ForceText=$false},
    @{Cell='E13'; Value='  -1.23%  '; ForceText=$false},
    @{Cell='D14'; Value='25.87'; ForceText=$false},
    @{Cell='E14'; Value='  -0.86%  '; ForceText=$false},
    @{Cell='D15'; Value='2.861.77'; ForceText=$false},
    @{Cell='E15'; Value='  -1.08%  '; ForceText=$false},
    @{Cell='D16'; Value='0.0000169'; ForceText=$false},
    @{Cell='E16'; Value='  -2.39%  '; ForceText=$false},
    @{Cell='D17'; Value='61.020.46'; ForceText=$false},
    @{Cell='E17'; Value='  -1.34%  '; ForceText=$false},
    @{Cell='D18'; Value='2.408.62'; ForceText=$false},
    @{Cell='E18'; Value='  -1.40%  '; ForceText=$false},
    @{Cell='D19'; Value='8.14'; ForceText=$false},
    @{Cell='E19'; Value='  +12.93%  '; ForceText=$false},
    @{Cell='D20'; Value='10.51'; ForceText=$false},
    @{Cell='E20'; Value='  -1.03%  '; ForceText=$false},
    @{Cell='D21'; Value='322.25'; ForceText=$false},
    @{Cell='E21'; Value='  -0.88%  '; ForceText=$false},
    @{Cell='E22'; Value='  -0.73%  '; ForceText=$false},
    @{Cell='D23'; Value='6.19'; ForceText=$false},
    @{Cell='E23'; Value='  +3.78%  '; ForceText=$false},
    @{Cell='E24'; Value='  -0.06%  '; ForceText=$false},
    @{Cell='E25'; Value='  -4.06%  '; ForceText=$false},
    @{Cell='D26'; Value='64.25'; ForceText=$false},
    @{Cell='E26'; Value='  -1.41%  '; ForceText=$false},
    @{Cell='D27'; Value='576.72'; ForceText=$false},
    @{Cell='E27'; Value='  -2.46%  '; ForceText=$false},
    @{Cell='E28'; Value='  -9.73%  '; ForceText=$false},
    @{Cell='D29'; Value='2.542.33'; ForceText=$false},
    @{Cell='E29'; Value='  -0.76%  '; ForceText=$false},
    @{Cell='D30'; Value='0.0₃0918'; ForceText=$false},
    @{Cell='E30'; Value='  -3.31%  '; ForceText=$false},
    @{Cell='D31'; Value='7.88'; ForceText=$false},
    @{Cell='E31'; Value='  -0.05%  '; ForceText=$false},
    @{Cell='E32'; Value='  -4.61%  '; ForceText=$false},
    @{Cell='E33'; Value='  -3.81%  '; ForceText=$false},
    @{Cell='E34'; Value='  -1.11%  '; ForceText=$false},
    @{Cell='E35'; Value='  +0.25%  '; ForceText=$false},
    @{Cell='D40'; Value='18.19'; ForceText=$false},
    @{Cell='E40'; Value='  -0.82%  '; ForceText=$false},
    @{Cell='D41'; Value='5.10'; ForceText=$true},
    @{Cell='E41'; Value='  -2.33%  '; ForceText=$false},
    @{Cell='E42'; Value='  +0.05%  '; ForceText=$false},
    @{Cell='D45'; Value='2.33'; ForceText=$false},
    @{Cell='E45'; Value='  -3.57%  '; ForceText=$false},
    @{Cell='D46'; Value='0.0₆0291'; ForceText=$false},
    @{Cell='E46'; Value='  +7.17%  '; ForceText=$false},
    @{Cell='D47'; Value='142.55'; ForceText=$false},
    @{Cell='E47'; Value='  +1.08%  '; ForceText=$false},
    @{Cell='D49'; Value='0.585'; ForceText=$false},
    @{Cell='E49'; Value='  -2.14%  '; ForceText=$false}
)

foreach ($u in $updates) {
    $r = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $r.NumberFormat = "@"
        $r.Value = $u.Value
        $r.ClearFormats()
    } else {
        $r.Value = $u.Value
    }
}

# Rows whose Coin/Link/Price/Volume were fully replaced (re-ranked rows)
$rowUpdates = @(
    @{Row=36; B='Monero'; C='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D='152.40'; DForceText=$true; E='  +0.00%  '},
    @{Row=37; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='1.40'; DForceText=$true; E='  -0.33%  '},
    @{Row=38; B='PolygonEcosystemToken'; C='https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'; D='0.367'; DForceText=$false; E='  -1.98%  '},
    @{Row=39; B='NEARProtocol'; C='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D='4.55'; DForceText=$false; E='  -5.03%  '},
    @{Row=43; B='OKB'; C='https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; D='41.20'; DForceText=$true; E='  -4.24%  '},
    @{Row=44; B='Stacks'; C='https://coinranking.com/coin/mMPrMcB7+stacks-stx'; D='1.65'; DForceText=$false; E='  -1.93%  '},
    @{Row=50; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.0499'; DForceText=$false; E='  -2.69%  '},
    @{Row=51; B='InjectiveProtocol'; C='https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; D='19.10'; DForceText=$true; E='  -2.61%  '}
)

foreach ($row in $rowUpdates) {
    $ws.Cells.Item($row.Row, 2).Value = $row.B
    $ws.Cells.Item($row.Row, 3).Value = $row.C
    $dCell = $ws.Cells.Item($row.Row, 4)
    if ($row.DForceText) {
        $dCell.NumberFormat = "@"
        $dCell.Value = $row.D
        $dCell.ClearFormats()
    } else {
        $dCell.Value = $row.D
    }
    $ws.Cells.Item($row.Row, 5).Value = $row.E
}

Write-Host "Applied $($updates.Count) cell updates and $($rowUpdates.Count) row updates"
